$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164; this shifts the existing rows 164-232
# down to 165-233 (preserving all their data/styles).
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with its data. Columns A, B, C, E,
# F, G, H, I, N, Q, R keep the same values as the template row (the row
# that used to be at 164, now at 165); columns D, J, K, L, M, O, P get
# new values.
$ws.Cells.Item(164, 1).Value = 3
$ws.Cells.Item(164, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(164, 3).Value = "Coquimbo"
$ws.Cells.Item(164, 4).Value = 44900
$ws.Cells.Item(164, 5).Value = 5
$ws.Cells.Item(164, 6).Value = 100112026
$ws.Cells.Item(164, 7).Value = "Haba"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 85
$ws.Cells.Item(164, 11).Value = 8000
$ws.Cells.Item(164, 12).Value = 8500
$ws.Cells.Item(164, 13).Value = 8235
$ws.Cells.Item(164, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(164, 15).Value = "Región Metropolitana"
$ws.Cells.Item(164, 16).Value = 329
$ws.Cells.Item(164, 17).Value = 25
$ws.Cells.Item(164, 18).Value = "Hortaliza"

# Keep the date-format (same numeric format used by the rest of column D)
# on the new row's date cell.
$ws.Cells.Item(164, 4).NumberFormat = $ws.Cells.Item(165, 4).NumberFormat()
